$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at position 259, pushing the existing rows
# 259..302 down to 263..306 (matches the new dimension A1:T306).
$ws.Rows(259).Insert()
$ws.Rows(259).Insert()
$ws.Rows(259).Insert()
$ws.Rows(259).Insert()

# Shared values for every row in this data block.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102004
$categoria = "Mandarina"
$unidad = "`$/bandeja 10 kilos"
$origen = "Provincia de Limarí"
$kgUnidad = 10

function Set-DataRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

Set-DataRow 259 44449 "Murcott" "Especial" 400 5000 5500 5250 525
Set-DataRow 260 44449 "Murcott" "Primera"  500 4000 4500 4250 425
Set-DataRow 261 44449 "Murcott" "Segunda"  400 3000 3500 3250 325
Set-DataRow 262 44449 "Murcott" "Tercera"  300 2000 2500 2250 225

Write-Output "inserted 4 rows at 259"
